$d = $word.ActiveDocument

# 1. Merge "gwynne.cs.ualberta.ca:" + "1521:CRS" + " will lead..." into one run (no text change)
$d.Content.Find.Execute(
    "gwynne.cs.ualberta.ca:1521:CRS will lead to the database stored at the given host ).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "gwynne.cs.ualberta.ca:1521:CRS will lead to the database stored at the given host ).",
    2)

# 2. Merge " " + "ARS's main interface calls..." into one run (no text change)
$d.Content.Find.Execute(
    " ARS" + [char]0x2019 + "s main interface calls each of the five major programs detailed below.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " ARS" + [char]0x2019 + "s main interface calls each of the five major programs detailed below.",
    2)

# 3. Merge APP1 source-code sentence runs (no text change)
$d.Content.Find.Execute(
    "APP1" + [char]0x2019 + "s source code is stored in ./apps/new_vehicle_registration.py",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "APP1" + [char]0x2019 + "s source code is stored in ./apps/new_vehicle_registration.py",
    2)

# 4. Merge APP2 "APP2's source code is stored " + "in ." + "/apps/auto_transaction.py" (no text change)
$d.Content.Find.Execute(
    "APP2" + [char]0x2019 + "s source code is stored in ./apps/auto_transaction.py",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "APP2" + [char]0x2019 + "s source code is stored in ./apps/auto_transaction.py",
    2)

# 5. "In APP4 ... appropriate." -> "In APP3 ... appropriate. There is also a small
#    conditions widget ... database." with the _GoBack bookmark repositioned right
#    after "In APP3".
$d.Content.Find.Execute(
    "In APP4 the user can enter the information to create a new license. The app ensures the unique constraint on SIN and License # are maintained when submitting data. This app also makes use of the",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In APP3 the user can enter the information to create a new license. The app ensures the unique constraint on SIN and License # are maintained when submitting data. This app also makes use of the",
    2)

$newTail = " There is also a small conditions widget to add new conditions to the database for use with the app, and the " + [char]0x201C + "?" + [char]0x201D + " button will let you see all conditions in the database."
$d.Content.Find.Execute(
    "You can also open the photo file you have selected and make sure that the picture is appropriate.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "You can also open the photo file you have selected and make sure that the picture is appropriate." + $newTail,
    2)

# Reposition the _GoBack bookmark to sit right after "In APP3" (collapsed range).
$bmRange = $d.Content
$bmRange.Find.Execute("In APP3")
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 6. Merge APP3 source-code sentence runs (no text change)
$d.Content.Find.Execute(
    "APP3" + [char]0x2019 + "s source code is stored in ./apps/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "APP3" + [char]0x2019 + "s source code is stored in ./apps/",
    2)

# 7. Violation-record paragraph: merge "vType" spell-check runs and shift the
#    "violation without a description..." sentence so it starts right after
#    "submit a " (same overall text, run boundaries change).
$q = [char]0x201C
$qe = [char]0x201D
$old7 = "The user can use APP4 to enter a Violation that has been issued. The format for date is preloaded and pressing the " + $q + "?" + $qe + " beside the " + $q + "Date Issued" + $qe + " column will allow the user to set the entry to the system time. Pressing the " + $q + "?" + $qe + " by the vType entry will allow the user to pull up a list of all the types of violations and their associated fines. You are allowed to submit a violation without a description and without a violator SIN. If you choose to submit without a violator SIN, the application will place the ticket on the primary owner of the entered VIN."
$new7 = "The user can use APP4 to enter a Violation that has been issued. The format for date is preloaded and pressing the " + $q + "?" + $qe + " beside the " + $q + "Date Issued" + $qe + " column will allow the user to set the entry to the system time. Pressing the " + $q + "?" + $qe + " by the vType entry will allow the user to pull up a list of all the types of violations and their associated fines. You are allowed to submit a violation without a description and without a violator SIN. If you choose to submit without a violator SIN, the application will place the ticket on the primary owner of the entered VIN."
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)

# 8. Merge APP4 source-code sentence runs (no text change)
$d.Content.Find.Execute(
    "APP4" + [char]0x2019 + "s source code is stored in ./apps/violation_records.py",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "APP4" + [char]0x2019 + "s source code is stored in ./apps/violation_records.py",
    2)

# 9. Merge "Searching for personal information ... etc)" runs (no text change)
$d.Content.Find.Execute(
    "Searching for personal information (i.e. Address, Birthday, driving class, etc)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Searching for personal information (i.e. Address, Birthday, driving class, etc)",
    2)

# 10. Merge "Searching for violation history ... etc) by searching..." runs (no text change)
$d.Content.Find.Execute(
    "Searching for violation history (i.e. Ticket Number, Issuing Officer, location, fine, etc) by searching on a SIN or license number.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Searching for violation history (i.e. Ticket Number, Issuing Officer, location, fine, etc) by searching on a SIN or license number.",
    2)

# 11. Merge APP5 source-code sentence runs (no text change)
$d.Content.Find.Execute(
    "APP5" + [char]0x2019 + "s source code is stored in ./apps/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "APP5" + [char]0x2019 + "s source code is stored in ./apps/",
    2)

# 12. Merge "3a i Design Choices" runs (no text change)
$d.Content.Find.Execute(
    "3a i Design Choices",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3a i Design Choices",
    2)

# 13. Merge "3a ii Git Hub" runs (no text change)
$d.Content.Find.Execute(
    "3a ii Git Hub",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3a ii Git Hub",
    2)
